$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 37, shifting existing rows 37..124 down to 38..125
$ws.Rows(37).Insert()

# Populate the new row 37 with the new record's data.
# Columns A,B,C,E,F,G,H,I,N,O,Q,R are constant across all rows in this table.
$ws.Range("A37").Value = 8
$ws.Range("B37").Value = "Terminal La Palmera de La Serena"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 44498
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = 100112037
$ws.Range("G37").Value = "Cebollín"
$ws.Range("H37").Value = "Sin especificar"
$ws.Range("I37").Value = "Primera"
$ws.Range("J37").Value = 3000
$ws.Range("K37").Value = 900
$ws.Range("L37").Value = 1000
$ws.Range("M37").Value = 950
$ws.Range("N37").Value = "$/paquete 6 unidades"
$ws.Range("O37").Value = "Provincia del Elquí"
$ws.Range("P37").Value = 158
$ws.Range("Q37").Value = 6
$ws.Range("R37").Value = "Hortaliza"
